# Populate the "Data" sheet / Table1 with the parking-lot column headers
# and a first data row, then grow the table to cover the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow Table1 (and its autofilter) from A1:E1 to A1:E2 so the new data
# row becomes part of the table.
$lo.Resize($ws.Range("A1:E2"))

# Header row (also renames the table's columns from the default "None").
$ws.Range("A1").Value = "Lot"
$ws.Range("B1").Value = "Status"
$ws.Range("C1").Value = "Date&Time"
$ws.Range("D1").Value = "BusySpaces"
$ws.Range("E1").Value = "FreeSpaces"

# First data row.
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 25
